{"js": "// Replace each two-digit-division expression with its updated value.\n// Every \"before\" value occurs exactly once in the document, so a direct\n// search + replace per pair (order independent) reproduces the diff.\nconst replacements = [\n  [\"65\u00f73=\", \"21\u00f79=\"],\n  [\"52\u00f74=\", \"30\u00f74=\"],\n  [\"50\u00f77=\", \"92\u00f76=\"],\n  [\"98\u00f74=\", \"47\u00f74=\"],\n  [\"58\u00f72=\", \"34\u00f75=\"],\n  [\"61\u00f72=\", \"44\u00f72=\"],\n  [\"32\u00f75=\", \"18\u00f74=\"],\n  [\"81\u00f73=\", \"53\u00f78=\"],\n  [\"32\u00f72=\", \"18\u00f76=\"],\n  [\"93\u00f73=\", \"52\u00f75=\"],\n  [\"31\u00f72=\", \"30\u00f78=\"],\n  [\"72\u00f74=\", \"26\u00f75=\"],\n  [\"41\u00f73=\", \"22\u00f72=\"],\n  [\"19\u00f74=\", \"17\u00f79=\"],\n  [\"18\u00f79=\", \"88\u00f78=\"],\n  [\"24\u00f74=\", \"52\u00f79=\"],\n  [\"35\u00f76=\", \"35\u00f75=\"],\n  [\"66\u00f78=\", \"28\u00f78=\"],\n  [\"73\u00f75=\", \"10\u00f75=\"],\n  [\"39\u00f73=\", \"59\u00f72=\"],\n  [\"48\u00f77=\", \"55\u00f77=\"],\n  [\"33\u00f73=\", \"96\u00f79=\"],\n  [\"37\u00f79=\", \"46\u00f73=\"],\n  [\"54\u00f72=\", \"66\u00f74=\"],\n  [\"57\u00f73=\", \"34\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-division expression with its updated value.\n# Every \"before\" string is unique within the document body, so a direct\n# Find/Replace per pair (order independent) reproduces the diff exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"65\u00f73=\", \"21\u00f79=\"),\n    @(\"52\u00f74=\", \"30\u00f74=\"),\n    @(\"50\u00f77=\", \"92\u00f76=\"),\n    @(\"98\u00f74=\", \"47\u00f74=\"),\n    @(\"58\u00f72=\", \"34\u00f75=\"),\n    @(\"61\u00f72=\", \"44\u00f72=\"),\n    @(\"32\u00f75=\", \"18\u00f74=\"),\n    @(\"81\u00f73=\", \"53\u00f78=\"),\n    @(\"32\u00f72=\", \"18\u00f76=\"),\n    @(\"93\u00f73=\", \"52\u00f75=\"),\n    @(\"31\u00f72=\", \"30\u00f78=\"),\n    @(\"72\u00f74=\", \"26\u00f75=\"),\n    @(\"41\u00f73=\", \"22\u00f72=\"),\n    @(\"19\u00f74=\", \"17\u00f79=\"),\n    @(\"18\u00f79=\", \"88\u00f78=\"),\n    @(\"24\u00f74=\", \"52\u00f79=\"),\n    @(\"35\u00f76=\", \"35\u00f75=\"),\n    @(\"66\u00f78=\", \"28\u00f78=\"),\n    @(\"73\u00f75=\", \"10\u00f75=\"),\n    @(\"39\u00f73=\", \"59\u00f72=\"),\n    @(\"48\u00f77=\", \"55\u00f77=\"),\n    @(\"33\u00f73=\", \"96\u00f79=\"),\n    @(\"37\u00f79=\", \"46\u00f73=\"),\n    @(\"54\u00f72=\", \"66\u00f74=\"),\n    @(\"57\u00f73=\", \"34\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        Write-Output (\"No match found for pair index \" + $replacements.IndexOf($pair))\n    }\n}\n\n"}
